# The deck's live theme part (ppt/theme/theme2.xml, "Integral" / "Red Violet")
# and its previously-unused twin (ppt/theme/theme1.xml, "Office Theme" / "Office")
# traded places: theme2.xml now carries the Office color scheme, theme1.xml now
# carries the Integral ("Red Violet") color scheme. The font scheme and format
# scheme are identical between the two themes, so only the 12-slot colour
# scheme actually changes in appearance.
#
# Apply the Office-theme colour values to the presentation's live theme
# (reached through any slide's ThemeColorScheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
